$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$lastRow = $tbl.Rows.Count

function Insert-CellRun($Cell, $ParaAttrs, $Text) {
    $cellRange = $Cell.Range
    $paraRange = $cellRange.Paragraphs.Item(1).Range
    $ins = $paraRange.Duplicate
    $ins.Collapse(1)  # wdCollapseStart

    $xml = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' +
           '<w:p' + $ParaAttrs + '>' +
           '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>' + $Text + '</w:t></w:r>' +
           '</w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$ins.InsertXML($xml)
}

Insert-CellRun $tbl.Cell($lastRow, 1) ' w14:paraId="05D3FB64" w14:textId="77777777" w:rsidR="00DA111A" w:rsidRPr="00DE7494" w:rsidRDefault="00DA111A"' "Jing Zhang"

Insert-CellRun $tbl.Cell($lastRow, 2) ' w14:paraId="305FDDAE" w14:textId="77777777" w:rsidR="00DA111A" w:rsidRDefault="00DA111A"' "Trends in the Use of Augmented Reality, Virtual Reality, and Mixed Reality in Surgical Research"

Insert-CellRun $tbl.Cell($lastRow, 3) ' w14:paraId="68C83872" w14:textId="77777777" w:rsidR="00DA111A" w:rsidRDefault="00DA111A"' "Global Bibliometric and Visualized Analysis"
